# Automatische test-sync: 2025-08-04 20:42:50
# Append a new test-mail log row to the "Logs" sheet and bump the matching
# category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 20

$logs.Range("A" + $newRow).Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Range("B" + $newRow).Value = "mailmind.test@zohomail.eu"
$logs.Range("C" + $newRow).Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Range("D" + $newRow).Value = "Inkoop / Bestellingen"
$logs.Range("E" + $newRow).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Range("F" + $newRow).Value = "2025-08-04 20:42:45"
$logs.Range("G" + $newRow).Value = "Ja"
$logs.Range("H" + $newRow).Value = "Ja"
$logs.Range("I" + $newRow).Value = "Nee"
$logs.Range("J" + $newRow).Value = "Nee"

# The conditional-formatting rules for D/G/H/I/J previously covered rows
# 2:19; extend each to include the freshly appended row 20.
$logs.Range("D2:D19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))
$logs.Range("G2:G19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))
$logs.Range("H2:H19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H20"))
$logs.Range("I2:I19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I20"))
$logs.Range("J2:J19").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J20"))

# Bump the Dashboard "Inkoop / Bestellingen" counter from 5 to 6.
$current = $dashboard.Range("B3").Value2
$dashboard.Range("B3").Value = $current + 1
